$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update electricity price (column B) for Namibia and Other rows
$ws.Range("B2").Value = 0.108
$ws.Range("B3").Value = 0.108

# Update interest rate columns (Solar/Wind/Plant/Infrastructure) to the new flat 6% rate
$ws.Range("D2").Value = 0.06
$ws.Range("F2").Value = 0.06
$ws.Range("H2").Value = 0.06
$ws.Range("J2").Value = 0.06

$ws.Range("D3").Value = 0.06
$ws.Range("F3").Value = 0.06
$ws.Range("H3").Value = 0.06
$ws.Range("J3").Value = 0.06

# Match the new active selection left by the author's session
[void]$ws.Range("A3").Select()
